$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the date number format is applied to the full Interval Date column (A2:A97)
$ws.Range("A2:A97").NumberFormat = "dd.mm.yyyy"

$ws.Cells.Item(2, 1).Value = 45327
$ws.Cells.Item(2, 2).Value = 13
$ws.Cells.Item(2, 3).Value = 3.371926546096802
$ws.Cells.Item(3, 1).Value = 45327
$ws.Cells.Item(3, 2).Value = 14
$ws.Cells.Item(3, 3).Value = 3.772702217102051
$ws.Cells.Item(4, 1).Value = 45327
$ws.Cells.Item(4, 2).Value = 15
$ws.Cells.Item(4, 3).Value = 3.558686494827271
$ws.Cells.Item(5, 1).Value = 45327
$ws.Cells.Item(5, 2).Value = 16
$ws.Cells.Item(5, 3).Value = 1.746927738189697
$ws.Cells.Item(6, 1).Value = 45327
$ws.Cells.Item(6, 2).Value = 17
$ws.Cells.Item(6, 3).Value = 0.02620400488376617
$ws.Cells.Item(7, 1).Value = 45327
$ws.Cells.Item(7, 2).Value = 18
$ws.Cells.Item(7, 3).Value = 0.02413088455796242
$ws.Cells.Item(8, 1).Value = 45327
$ws.Cells.Item(8, 2).Value = 19
$ws.Cells.Item(8, 3).Value = 0.02413088455796242
$ws.Cells.Item(9, 1).Value = 45327
$ws.Cells.Item(9, 2).Value = 20
$ws.Cells.Item(9, 3).Value = 0.02413088455796242
$ws.Cells.Item(10, 1).Value = 45327
$ws.Cells.Item(10, 2).Value = 21
$ws.Cells.Item(10, 3).Value = 0.02413088455796242
$ws.Cells.Item(11, 1).Value = 45327
$ws.Cells.Item(11, 2).Value = 22
$ws.Cells.Item(11, 3).Value = 0.02413088455796242
$ws.Cells.Item(12, 1).Value = 45327
$ws.Cells.Item(12, 2).Value = 23
$ws.Cells.Item(12, 3).Value = 0.02413088455796242
$ws.Cells.Item(13, 1).Value = 45328
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 0.02396511659026146
$ws.Cells.Item(14, 1).Value = 45328
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Value = 0.02396511659026146
$ws.Cells.Item(15, 1).Value = 45328
$ws.Cells.Item(15, 2).Value = 2
$ws.Cells.Item(15, 3).Value = 0.02396511659026146
$ws.Cells.Item(16, 1).Value = 45328
$ws.Cells.Item(16, 2).Value = 3
$ws.Cells.Item(16, 3).Value = 0.02396511659026146
$ws.Cells.Item(17, 1).Value = 45328
$ws.Cells.Item(17, 2).Value = 4
$ws.Cells.Item(17, 3).Value = 0.02396511659026146
$ws.Cells.Item(18, 1).Value = 45328
$ws.Cells.Item(18, 2).Value = 5
$ws.Cells.Item(18, 3).Value = 0.02396511659026146
$ws.Cells.Item(19, 1).Value = 45328
$ws.Cells.Item(19, 2).Value = 6
$ws.Cells.Item(19, 3).Value = 0.02396511659026146
$ws.Cells.Item(20, 1).Value = 45328
$ws.Cells.Item(20, 2).Value = 7
$ws.Cells.Item(20, 3).Value = 0.0250907875597477
$ws.Cells.Item(21, 1).Value = 45328
$ws.Cells.Item(21, 2).Value = 8
$ws.Cells.Item(21, 3).Value = 0.8482181429862976
$ws.Cells.Item(22, 1).Value = 45328
$ws.Cells.Item(22, 2).Value = 9
$ws.Cells.Item(22, 3).Value = 2.440028667449951
$ws.Cells.Item(23, 1).Value = 45328
$ws.Cells.Item(23, 2).Value = 10
$ws.Cells.Item(23, 3).Value = 3.281532526016235
$ws.Cells.Item(24, 1).Value = 45328
$ws.Cells.Item(24, 2).Value = 11
$ws.Cells.Item(24, 3).Value = 3.795012235641479
$ws.Cells.Item(25, 1).Value = 45328
$ws.Cells.Item(25, 2).Value = 12
$ws.Cells.Item(25, 3).Value = 3.81330680847168
$ws.Cells.Item(26, 1).Value = 45328
$ws.Cells.Item(26, 2).Value = 13
$ws.Cells.Item(26, 3).Value = 3.16619086265564
$ws.Cells.Item(27, 1).Value = 45328
$ws.Cells.Item(27, 2).Value = 14
$ws.Cells.Item(27, 3).Value = 3.527982950210571
$ws.Cells.Item(28, 1).Value = 45328
$ws.Cells.Item(28, 2).Value = 15
$ws.Cells.Item(28, 3).Value = 3.410948991775513
$ws.Cells.Item(29, 1).Value = 45328
$ws.Cells.Item(29, 2).Value = 16
$ws.Cells.Item(29, 3).Value = 2.021544456481934
$ws.Cells.Item(30, 1).Value = 45328
$ws.Cells.Item(30, 2).Value = 17
$ws.Cells.Item(30, 3).Value = 0.0350673496723175
$ws.Cells.Item(31, 1).Value = 45328
$ws.Cells.Item(31, 2).Value = 18
$ws.Cells.Item(31, 3).Value = 0.02413088455796242
$ws.Cells.Item(32, 1).Value = 45328
$ws.Cells.Item(32, 2).Value = 19
$ws.Cells.Item(32, 3).Value = 0.02347593382000923
$ws.Cells.Item(33, 1).Value = 45328
$ws.Cells.Item(33, 2).Value = 20
$ws.Cells.Item(33, 3).Value = 0.02347593382000923
$ws.Cells.Item(34, 1).Value = 45328
$ws.Cells.Item(34, 2).Value = 21
$ws.Cells.Item(34, 3).Value = 0.02347593382000923
$ws.Cells.Item(35, 1).Value = 45328
$ws.Cells.Item(35, 2).Value = 22
$ws.Cells.Item(35, 3).Value = 0.02320006862282753
$ws.Cells.Item(36, 1).Value = 45328
$ws.Cells.Item(36, 2).Value = 23
$ws.Cells.Item(36, 3).Value = 0.02320006862282753
$ws.Cells.Item(37, 1).Value = 45329
$ws.Cells.Item(37, 2).Value = 0
$ws.Cells.Item(37, 3).Value = 0.02331016585230827
$ws.Cells.Item(38, 1).Value = 45329
$ws.Cells.Item(38, 2).Value = 1
$ws.Cells.Item(38, 3).Value = 0.02396511659026146
$ws.Cells.Item(39, 1).Value = 45329
$ws.Cells.Item(39, 2).Value = 2
$ws.Cells.Item(39, 3).Value = 0.02396511659026146
$ws.Cells.Item(40, 1).Value = 45329
$ws.Cells.Item(40, 2).Value = 3
$ws.Cells.Item(40, 3).Value = 0.02396511659026146
$ws.Cells.Item(41, 1).Value = 45329
$ws.Cells.Item(41, 2).Value = 4
$ws.Cells.Item(41, 3).Value = 0.02396511659026146
$ws.Cells.Item(42, 1).Value = 45329
$ws.Cells.Item(42, 2).Value = 5
$ws.Cells.Item(42, 3).Value = 0.02396511659026146
$ws.Cells.Item(43, 1).Value = 45329
$ws.Cells.Item(43, 2).Value = 6
$ws.Cells.Item(43, 3).Value = 0.02396511659026146
$ws.Cells.Item(44, 1).Value = 45329
$ws.Cells.Item(44, 2).Value = 7
$ws.Cells.Item(44, 3).Value = 0.02602160349488258
$ws.Cells.Item(45, 1).Value = 45329
$ws.Cells.Item(45, 2).Value = 8
$ws.Cells.Item(45, 3).Value = 1.76830267906189
$ws.Cells.Item(46, 1).Value = 45329
$ws.Cells.Item(46, 2).Value = 9
$ws.Cells.Item(46, 3).Value = 3.675012588500977
$ws.Cells.Item(47, 1).Value = 45329
$ws.Cells.Item(47, 2).Value = 10
$ws.Cells.Item(47, 3).Value = 3.513988733291626
$ws.Cells.Item(48, 1).Value = 45329
$ws.Cells.Item(48, 2).Value = 11
$ws.Cells.Item(48, 3).Value = 3.19330620765686
$ws.Cells.Item(49, 1).Value = 45329
$ws.Cells.Item(49, 2).Value = 12
$ws.Cells.Item(49, 3).Value = 3.115317821502686
$ws.Cells.Item(50, 1).Value = 45329
$ws.Cells.Item(50, 2).Value = 13
$ws.Cells.Item(50, 3).Value = 3.102758884429932
$ws.Cells.Item(51, 1).Value = 45329
$ws.Cells.Item(51, 2).Value = 14
$ws.Cells.Item(51, 3).Value = 3.241005897521973
$ws.Cells.Item(52, 1).Value = 45329
$ws.Cells.Item(52, 2).Value = 15
$ws.Cells.Item(52, 3).Value = 3.633407115936279
$ws.Cells.Item(53, 1).Value = 45329
$ws.Cells.Item(53, 2).Value = 16
$ws.Cells.Item(53, 3).Value = 1.129147887229919
$ws.Cells.Item(54, 1).Value = 45329
$ws.Cells.Item(54, 2).Value = 17
$ws.Cells.Item(54, 3).Value = 0.01550001557916403
$ws.Cells.Item(55, 1).Value = 45329
$ws.Cells.Item(55, 2).Value = 18
$ws.Cells.Item(55, 3).Value = 0.01349807716906071
$ws.Cells.Item(56, 1).Value = 45329
$ws.Cells.Item(56, 2).Value = 19
$ws.Cells.Item(56, 3).Value = 0.01207545306533575
$ws.Cells.Item(57, 1).Value = 45329
$ws.Cells.Item(57, 2).Value = 20
$ws.Cells.Item(57, 3).Value = 0.009953487664461136
$ws.Cells.Item(58, 1).Value = 45329
$ws.Cells.Item(58, 2).Value = 21
$ws.Cells.Item(58, 3).Value = 0.009953487664461136
$ws.Cells.Item(59, 1).Value = 45329
$ws.Cells.Item(59, 2).Value = 22
$ws.Cells.Item(59, 3).Value = 0.01113460119813681
$ws.Cells.Item(60, 1).Value = 45329
$ws.Cells.Item(60, 2).Value = 23
$ws.Cells.Item(60, 3).Value = 0.01207545306533575
$ws.Cells.Item(61, 1).Value = 45330
$ws.Cells.Item(61, 2).Value = 0
$ws.Cells.Item(61, 3).Value = 0.01280297338962555
$ws.Cells.Item(62, 1).Value = 45330
$ws.Cells.Item(62, 2).Value = 1
$ws.Cells.Item(62, 3).Value = 0.02396511659026146
$ws.Cells.Item(63, 1).Value = 45330
$ws.Cells.Item(63, 2).Value = 2
$ws.Cells.Item(63, 3).Value = 0.02396511659026146
$ws.Cells.Item(64, 1).Value = 45330
$ws.Cells.Item(64, 2).Value = 3
$ws.Cells.Item(64, 3).Value = 0.02396511659026146
$ws.Cells.Item(65, 1).Value = 45330
$ws.Cells.Item(65, 2).Value = 4
$ws.Cells.Item(65, 3).Value = 0.02396511659026146
$ws.Cells.Item(66, 1).Value = 45330
$ws.Cells.Item(66, 2).Value = 5
$ws.Cells.Item(66, 3).Value = 0.02396511659026146
$ws.Cells.Item(67, 1).Value = 45330
$ws.Cells.Item(67, 2).Value = 6
$ws.Cells.Item(67, 3).Value = 0.02396511659026146
$ws.Cells.Item(68, 1).Value = 45330
$ws.Cells.Item(68, 2).Value = 7
$ws.Cells.Item(68, 3).Value = 0.02602160349488258
$ws.Cells.Item(69, 1).Value = 45330
$ws.Cells.Item(69, 2).Value = 8
$ws.Cells.Item(69, 3).Value = 2.102064847946167
$ws.Cells.Item(70, 1).Value = 45330
$ws.Cells.Item(70, 2).Value = 9
$ws.Cells.Item(70, 3).Value = 3.77857232093811
$ws.Cells.Item(71, 1).Value = 45330
$ws.Cells.Item(71, 2).Value = 10
$ws.Cells.Item(71, 3).Value = 3.54853367805481
$ws.Cells.Item(72, 1).Value = 45330
$ws.Cells.Item(72, 2).Value = 11
$ws.Cells.Item(72, 3).Value = 3.367605686187744
$ws.Cells.Item(73, 1).Value = 45330
$ws.Cells.Item(73, 2).Value = 12
$ws.Cells.Item(73, 3).Value = 3.180029630661011
$ws.Cells.Item(74, 1).Value = 45330
$ws.Cells.Item(74, 2).Value = 13
$ws.Cells.Item(74, 3).Value = 3.318439960479736
$ws.Cells.Item(75, 1).Value = 45330
$ws.Cells.Item(75, 2).Value = 14
$ws.Cells.Item(75, 3).Value = 3.821091890335083
$ws.Cells.Item(76, 1).Value = 45330
$ws.Cells.Item(76, 2).Value = 15
$ws.Cells.Item(76, 3).Value = 3.424803495407104
$ws.Cells.Item(77, 1).Value = 45330
$ws.Cells.Item(77, 2).Value = 16
$ws.Cells.Item(77, 3).Value = 1.72944450378418
$ws.Cells.Item(78, 1).Value = 45330
$ws.Cells.Item(78, 2).Value = 17
$ws.Cells.Item(78, 3).Value = 0.02620400488376617
$ws.Cells.Item(79, 1).Value = 45330
$ws.Cells.Item(79, 2).Value = 18
$ws.Cells.Item(79, 3).Value = 0.02413088455796242
$ws.Cells.Item(80, 1).Value = 45330
$ws.Cells.Item(80, 2).Value = 19
$ws.Cells.Item(80, 3).Value = 0.02413088455796242
$ws.Cells.Item(81, 1).Value = 45330
$ws.Cells.Item(81, 2).Value = 20
$ws.Cells.Item(81, 3).Value = 0.02413088455796242
$ws.Cells.Item(82, 1).Value = 45330
$ws.Cells.Item(82, 2).Value = 21
$ws.Cells.Item(82, 3).Value = 0.02413088455796242
$ws.Cells.Item(83, 1).Value = 45330
$ws.Cells.Item(83, 2).Value = 22
$ws.Cells.Item(83, 3).Value = 0.02413088455796242
$ws.Cells.Item(84, 1).Value = 45330
$ws.Cells.Item(84, 2).Value = 23
$ws.Cells.Item(84, 3).Value = 0.02413088455796242
$ws.Cells.Item(85, 1).Value = 45331
$ws.Cells.Item(85, 2).Value = 0
$ws.Cells.Item(85, 3).Value = 3.394882440567017
$ws.Cells.Item(86, 1).Value = 45331
$ws.Cells.Item(86, 2).Value = 1
$ws.Cells.Item(86, 3).Value = 3.462519645690918
$ws.Cells.Item(87, 1).Value = 45331
$ws.Cells.Item(87, 2).Value = 2
$ws.Cells.Item(87, 3).Value = 3.462519645690918
$ws.Cells.Item(88, 1).Value = 45331
$ws.Cells.Item(88, 2).Value = 3
$ws.Cells.Item(88, 3).Value = 3.462519645690918
$ws.Cells.Item(89, 1).Value = 45331
$ws.Cells.Item(89, 2).Value = 4
$ws.Cells.Item(89, 3).Value = 3.462519645690918
$ws.Cells.Item(90, 1).Value = 45331
$ws.Cells.Item(90, 2).Value = 5
$ws.Cells.Item(90, 3).Value = 3.484997272491455
$ws.Cells.Item(91, 1).Value = 45331
$ws.Cells.Item(91, 2).Value = 6
$ws.Cells.Item(91, 3).Value = 3.462519645690918
$ws.Cells.Item(92, 1).Value = 45331
$ws.Cells.Item(92, 2).Value = 7
$ws.Cells.Item(92, 3).Value = 3.412805318832397
$ws.Cells.Item(93, 1).Value = 45331
$ws.Cells.Item(93, 2).Value = 8
$ws.Cells.Item(93, 3).Value = 3.383103132247925
$ws.Cells.Item(94, 1).Value = 45331
$ws.Cells.Item(94, 2).Value = 9
$ws.Cells.Item(94, 3).Value = 3.338872671127319
$ws.Cells.Item(95, 1).Value = 45331
$ws.Cells.Item(95, 2).Value = 10
$ws.Cells.Item(95, 3).Value = 3.318135023117065
$ws.Cells.Item(96, 1).Value = 45331
$ws.Cells.Item(96, 2).Value = 11
$ws.Cells.Item(96, 3).Value = 3.216139793395996
$ws.Cells.Item(97, 1).Value = 45331
$ws.Cells.Item(97, 2).Value = 12
$ws.Cells.Item(97, 3).Value = 3.102933168411255
